# Regenerate the "K" column (column G) values for the kelly_zack 2024
# save-data sheet. The pipeline that produces this workbook now computes
# K from the Strike# series (std/mean) instead of the previous formula,
# so the recalculated K values below replace the stale ones written to
# column G for rows 2-52 (row 1 is the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 4
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 3
    25 = 2
    26 = 2
    27 = 1
    28 = 3
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 1
    34 = 2
    35 = 3
    36 = 2
    37 = 3
    38 = 0
    39 = 2
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 2
    45 = 0
    46 = 2
    47 = 0
    48 = 2
    49 = 1
    50 = 3
    51 = 2
    52 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
